# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "75.581.72"
$ws.Range("E2").Value = "  +8.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.714.72"
$ws.Range("E3").Value = "  +11.64%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "187.16"
$ws.Range("E5").Value = "  +12.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "590.62"
$ws.Range("E6").Value = "  +4.59%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.541"
$ws.Range("E8").Value = "  +5.31%  "
$ws.Range("E9").Value = "  +15.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.713.31"
$ws.Range("E10").Value = "  +11.61%  "
$ws.Range("E11").Value = "  +1.32%  "
$ws.Range("E12").Value = "  +8.37%  "
$ws.Range("E13").Value = "  +1.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.214.59"
$ws.Range("E14").Value = "  +11.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "75.426.83"
$ws.Range("E15").Value = "  +8.37%  "
$ws.Range("E16").Value = "  +6.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.00"
$ws.Range("E17").Value = "  +12.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.718.58"
$ws.Range("E18").Value = "  +11.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.32"
$ws.Range("E19").Value = "  +29.48%  "
$ws.Range("E20").Value = "  +11.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "378.03"
$ws.Range("E21").Value = "  +9.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.31"
$ws.Range("E22").Value = "  +14.80%  "
$ws.Range("E23").Value = "  +6.31%  "
$ws.Range("E24").Value = "  +4.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.07"
$ws.Range("E25").Value = "  +7.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("E27").Value = "  +10.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.60"
$ws.Range("E28").Value = "  +13.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.855.75"
$ws.Range("E29").Value = "  +11.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.997"
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0989"
$ws.Range("E31").Value = "  +15.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "522.94"
$ws.Range("E32").Value = "  +14.31%  "
$ws.Range("E33").Value = "  +12.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.87"
$ws.Range("E34").Value = "  +6.74%  "
$ws.Range("E35").Value = "  +10.90%  "
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("E37").Value = "  +7.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.08"
$ws.Range("E38").Value = "  +1.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.58"
$ws.Range("E39").Value = "  +7.25%  "
$ws.Range("E40").Value = "  +1.30%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "173.10"
$ws.Range("E42").Value = "  +27.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.06"
$ws.Range("E43").Value = "  +14.42%  "
$ws.Range("E44").Value = "  +12.99%  "
$ws.Range("E45").Value = "  +9.85%  "
$ws.Range("E46").Value = "  +13.14%  "
$ws.Range("E47").Value = "  +14.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "39.20"
$ws.Range("E48").Value = "  +2.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0857"
$ws.Range("E49").Value = "  +18.49%  "
$ws.Range("E50").Value = "  +9.27%  "
$ws.Range("E51").Value = "  +11.67%  "
